# Auto-generated: updates currentAveragePrice / LevePrice / LeveProfit columns (H-N)
# across multiple worksheets, per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 3673.1667
$ws.Cells.Item(19, 9).Value = 4998
$ws.Cells.Item(19, 10).Value = 2017.125
$ws.Cells.Item(19, 11).Value = 4998
$ws.Cells.Item(19, 12).Value = 2017.125
$ws.Cells.Item(19, 13).Value = -4823
$ws.Cells.Item(19, 14).Value = -2367.125
$ws.Cells.Item(33, 8).Value = 339.0345
$ws.Cells.Item(33, 9).Value = 355.2963
$ws.Cells.Item(33, 11).Value = 355.2963
$ws.Cells.Item(33, 13).Value = -126.2963
$ws.Cells.Item(41, 8).Value = 2833.2173
$ws.Cells.Item(41, 9).Value = 75.666664
$ws.Cells.Item(41, 10).Value = 4605.9287
$ws.Cells.Item(41, 11).Value = 75.666664
$ws.Cells.Item(41, 12).Value = 4605.9287
$ws.Cells.Item(41, 13).Value = 364.333336
$ws.Cells.Item(41, 14).Value = -5485.9287
$ws.Cells.Item(98, 8).Value = 5713.8096
$ws.Cells.Item(98, 9).Value = 6421.6665
$ws.Cells.Item(98, 10).Value = 1466.6666
$ws.Cells.Item(98, 11).Value = 6421.6665
$ws.Cells.Item(98, 12).Value = 1466.6666
$ws.Cells.Item(98, 13).Value = -4923.6665
$ws.Cells.Item(98, 14).Value = -4462.6666
$ws.Cells.Item(122, 8).Value = 5713.8096
$ws.Cells.Item(122, 9).Value = 6421.6665
$ws.Cells.Item(122, 10).Value = 1466.6666
$ws.Cells.Item(122, 11).Value = 19264.9995
$ws.Cells.Item(122, 12).Value = 4399.9998
$ws.Cells.Item(122, 13).Value = -16814.9995
$ws.Cells.Item(122, 14).Value = -9299.9998
$ws.Cells.Item(129, 8).Value = 17715.283
$ws.Cells.Item(129, 9).Value = 524.5333000000001
$ws.Cells.Item(129, 11).Value = 1573.5999
$ws.Cells.Item(129, 13).Value = 3426.4001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3178.3462
$ws.Cells.Item(20, 9).Value = 3873
$ws.Cells.Item(20, 10).Value = 2367.9167
$ws.Cells.Item(20, 11).Value = 3873
$ws.Cells.Item(20, 12).Value = 2367.9167
$ws.Cells.Item(20, 13).Value = -3626
$ws.Cells.Item(20, 14).Value = -2861.9167
$ws.Cells.Item(80, 8).Value = 177.74074
$ws.Cells.Item(80, 9).Value = 250.42857
$ws.Cells.Item(80, 10).Value = 152.3
$ws.Cells.Item(80, 11).Value = 250.42857
$ws.Cells.Item(80, 12).Value = 152.3
$ws.Cells.Item(80, 13).Value = 747.57143
$ws.Cells.Item(80, 14).Value = -2148.3
$ws.Cells.Item(83, 8).Value = 177.74074
$ws.Cells.Item(83, 9).Value = 250.42857
$ws.Cells.Item(83, 10).Value = 152.3
$ws.Cells.Item(83, 11).Value = 1252.14285
$ws.Cells.Item(83, 12).Value = 761.5
$ws.Cells.Item(83, 13).Value = 3739.85715
$ws.Cells.Item(83, 14).Value = -10745.5
$ws.Cells.Item(94, 8).Value = 739.4167
$ws.Cells.Item(94, 9).Value = 724.8182
$ws.Cells.Item(94, 11).Value = 724.8182
$ws.Cells.Item(94, 13).Value = -273.8182
$ws.Cells.Item(134, 8).Value = 18009.28
$ws.Cells.Item(134, 9).Value = 24442.373
$ws.Cells.Item(134, 11).Value = 73327.11900000001
$ws.Cells.Item(134, 13).Value = -70792.11900000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1101.8334
$ws.Cells.Item(16, 9).Value = 1101.8334
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 1101.8334
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -814.8334
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(50, 8).Value = 12104.363
$ws.Cells.Item(50, 9).Value = 5000
$ws.Cells.Item(50, 10).Value = 12814.8
$ws.Cells.Item(50, 11).Value = 5000
$ws.Cells.Item(50, 12).Value = 12814.8
$ws.Cells.Item(50, 13).Value = -4375
$ws.Cells.Item(50, 14).Value = -14064.8
$ws.Cells.Item(51, 8).Value = 24495
$ws.Cells.Item(51, 10).Value = 24495
$ws.Cells.Item(51, 12).Value = 24495
$ws.Cells.Item(51, 14).Value = -25967
$ws.Cells.Item(60, 8).Value = 15339.182
$ws.Cells.Item(60, 10).Value = 15339.182
$ws.Cells.Item(60, 12).Value = 15339.182
$ws.Cells.Item(60, 14).Value = -16361.182
$ws.Cells.Item(61, 8).Value = 24495
$ws.Cells.Item(61, 10).Value = 24495
$ws.Cells.Item(61, 12).Value = 24495
$ws.Cells.Item(61, 14).Value = -25191
$ws.Cells.Item(86, 8).Value = 100003220
$ws.Cells.Item(86, 9).Value = 200001180
$ws.Cells.Item(86, 10).Value = 5259.6
$ws.Cells.Item(86, 11).Value = 200001180
$ws.Cells.Item(86, 12).Value = 5259.6
$ws.Cells.Item(86, 13).Value = -200000057
$ws.Cells.Item(86, 14).Value = -7505.6
$ws.Cells.Item(89, 8).Value = 100003220
$ws.Cells.Item(89, 9).Value = 200001180
$ws.Cells.Item(89, 10).Value = 5259.6
$ws.Cells.Item(89, 11).Value = 1000005900
$ws.Cells.Item(89, 12).Value = 26298
$ws.Cells.Item(89, 13).Value = -1000000284
$ws.Cells.Item(89, 14).Value = -37530
$ws.Cells.Item(94, 8).Value = 58823856
$ws.Cells.Item(94, 9).Value = 166666930
$ws.Cells.Item(94, 10).Value = 358.81818
$ws.Cells.Item(94, 11).Value = 166666930
$ws.Cells.Item(94, 12).Value = 358.81818
$ws.Cells.Item(94, 13).Value = -166666479
$ws.Cells.Item(94, 14).Value = -1260.81818
$ws.Cells.Item(99, 8).Value = 47521.09
$ws.Cells.Item(99, 10).Value = 2296.75
$ws.Cells.Item(99, 12).Value = 2296.75
$ws.Cells.Item(99, 14).Value = -5292.75
$ws.Cells.Item(105, 8).Value = 1307.8572
$ws.Cells.Item(105, 9).Value = 1049.875
$ws.Cells.Item(105, 11).Value = 1049.875
$ws.Cells.Item(105, 13).Value = 697.125
$ws.Cells.Item(107, 8).Value = 63276.625
$ws.Cells.Item(107, 9).Value = 77670.234
$ws.Cells.Item(107, 10).Value = 904.3333
$ws.Cells.Item(107, 11).Value = 77670.234
$ws.Cells.Item(107, 12).Value = 904.3333
$ws.Cells.Item(107, 13).Value = -75750.234
$ws.Cells.Item(107, 14).Value = -4744.3333
$ws.Cells.Item(113, 8).Value = 1101.8334
$ws.Cells.Item(113, 9).Value = 1101.8334
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 1101.8334
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = 1068.1666
$ws.Cells.Item(113, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 47521.09
$ws.Cells.Item(126, 10).Value = 2296.75
$ws.Cells.Item(126, 12).Value = 6890.25
$ws.Cells.Item(126, 14).Value = -11830.25
$ws.Cells.Item(134, 8).Value = 1260.2413
$ws.Cells.Item(134, 9).Value = 1237.5454
$ws.Cells.Item(134, 10).Value = 1331.5714
$ws.Cells.Item(134, 11).Value = 3712.6362
$ws.Cells.Item(134, 12).Value = 3994.7142
$ws.Cells.Item(134, 13).Value = -1177.6362
$ws.Cells.Item(134, 14).Value = -9064.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 1002.0909
$ws.Cells.Item(92, 9).Value = 1048
$ws.Cells.Item(92, 10).Value = 963.8333
$ws.Cells.Item(92, 11).Value = 3144
$ws.Cells.Item(92, 12).Value = 2891.4999
$ws.Cells.Item(92, 13).Value = -1896
$ws.Cells.Item(92, 14).Value = -5387.4999
$ws.Cells.Item(122, 8).Value = 2348.1
$ws.Cells.Item(122, 10).Value = 1668.1428
$ws.Cells.Item(122, 12).Value = 15013.2852
$ws.Cells.Item(122, 14).Value = -19913.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1542.7142
$ws.Cells.Item(102, 9).Value = 1359.8
$ws.Cells.Item(102, 10).Value = 2000
$ws.Cells.Item(102, 11).Value = 1359.8
$ws.Cells.Item(102, 12).Value = 2000
$ws.Cells.Item(102, 13).Value = 262.2
$ws.Cells.Item(102, 14).Value = -5244
$ws.Cells.Item(113, 8).Value = 31250994
$ws.Cells.Item(113, 10).Value = 1137
$ws.Cells.Item(113, 12).Value = 1137
$ws.Cells.Item(113, 14).Value = -5477

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 14).ClearContents()
$ws.Cells.Item(61, 8).Value = 37038780
$ws.Cells.Item(61, 9).Value = 1615.1666
$ws.Cells.Item(61, 10).Value = 111113110
$ws.Cells.Item(61, 11).Value = 1615.1666
$ws.Cells.Item(61, 12).Value = 111113110
$ws.Cells.Item(61, 13).Value = -1413.1666
$ws.Cells.Item(61, 14).Value = -111113514
$ws.Cells.Item(113, 8).Value = 37038780
$ws.Cells.Item(113, 9).Value = 1615.1666
$ws.Cells.Item(113, 10).Value = 111113110
$ws.Cells.Item(113, 11).Value = 1615.1666
$ws.Cells.Item(113, 12).Value = 111113110
$ws.Cells.Item(113, 13).Value = 554.8334
$ws.Cells.Item(113, 14).Value = -111117450

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2226.037
$ws.Cells.Item(122, 9).Value = 2179.7273
$ws.Cells.Item(122, 11).Value = 6539.1819
$ws.Cells.Item(122, 13).Value = -4089.1819
